$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1499.6666
$ws.Range("I18").Value = 1499.6666
$ws.Range("K18").Value = 1499.6666
$ws.Range("M18").Value = -1215.6666
$ws.Range("H32").Value = 1018
$ws.Range("J32").Value = 972.5
$ws.Range("L32").Value = 972.5
$ws.Range("N32").Value = -1624.5
$ws.Range("H38").Value = 414.75
$ws.Range("J38").Value = 998
$ws.Range("L38").Value = 2994
$ws.Range("N38").Value = -3738
$ws.Range("H112").Value = 79342
$ws.Range("I112").Value = 1091.6666
$ws.Range("K112").Value = 3274.9998
$ws.Range("M112").Value = -2166.9998
$ws.Range("H132").Value = 812.8387
$ws.Range("I132").Value = 822.76666
$ws.Range("K132").Value = 2468.29998
$ws.Range("M132").Value = 61.70002000000022
$ws.Range("H135").Value = 395.92
$ws.Range("I135").Value = 426.22726
$ws.Range("J135").Value = 173.66667
$ws.Range("K135").Value = 3836.04534
$ws.Range("L135").Value = 1563.00003
$ws.Range("M135").Value = -1301.04534
$ws.Range("N135").Value = -6633.00003
$ws.Range("H138").Value = 3682.09
$ws.Range("J138").Value = 3756.8645
$ws.Range("L138").Value = 11270.5935
$ws.Range("N138").Value = -21550.5935
$ws.Range("H141").Value = 5901.2666
$ws.Range("I141").Value = 5929.9287
$ws.Range("K141").Value = 17789.7861
$ws.Range("M141").Value = -12609.7861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 68996.75
$ws.Range("J44").Value = 68996.75
$ws.Range("L44").Value = 68996.75
$ws.Range("N44").Value = -69972.75
$ws.Range("H45").Value = 3598.4167
$ws.Range("I45").Value = 3608.5789
$ws.Range("J45").Value = 3559.8
$ws.Range("K45").Value = 3608.5789
$ws.Range("L45").Value = 3559.8
$ws.Range("M45").Value = -3231.5789
$ws.Range("N45").Value = -4313.8
$ws.Range("H55").Value = 50558.555
$ws.Range("J55").Value = 62996.5
$ws.Range("L55").Value = 62996.5
$ws.Range("N55").Value = -63626.5
$ws.Range("H61").Value = 2790.1052
$ws.Range("I61").Value = 2517
$ws.Range("J61").Value = 3999.5715
$ws.Range("K61").Value = 2517
$ws.Range("L61").Value = 3999.5715
$ws.Range("M61").Value = -2305
$ws.Range("N61").Value = -4423.5715
$ws.Range("H80").Value = 110016.336
$ws.Range("J80").Value = 109999.6
$ws.Range("L80").Value = 109999.6
$ws.Range("N80").Value = -111995.6
$ws.Range("H83").Value = 110016.336
$ws.Range("J83").Value = 109999.6
$ws.Range("L83").Value = 329998.8
$ws.Range("N83").Value = -339982.8
$ws.Range("H97").Value = 1097.9584
$ws.Range("J97").Value = 2550
$ws.Range("L97").Value = 2550
$ws.Range("N97").Value = -3542
$ws.Range("H136").Value = 2790.1052
$ws.Range("I136").Value = 2517
$ws.Range("J136").Value = 3999.5715
$ws.Range("K136").Value = 7551
$ws.Range("L136").Value = 11998.7145
$ws.Range("M136").Value = -5001
$ws.Range("N136").Value = -17098.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3177590.2
$ws.Range("I134").Value = 3923788.2
$ws.Range("K134").Value = 11771364.6
$ws.Range("M134").Value = -11768829.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1884.3721
$ws.Range("I31").Value = 2005.8334
$ws.Range("K31").Value = 2005.8334
$ws.Range("M31").Value = -1710.8334
$ws.Range("H34").Value = 1884.3721
$ws.Range("I34").Value = 2005.8334
$ws.Range("K34").Value = 2005.8334
$ws.Range("M34").Value = -1803.8334
$ws.Range("H58").Value = 4271.5884
$ws.Range("I58").Value = 3006.6
$ws.Range("J58").Value = 4798.6665
$ws.Range("K58").Value = 3006.6
$ws.Range("L58").Value = 4798.6665
$ws.Range("M58").Value = -2803.6
$ws.Range("N58").Value = -5204.6665
$ws.Range("H99").Value = 2999.6667
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 3999.5
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 3999.5
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -6995.5
$ws.Range("H104").Value = 64963
$ws.Range("J104").Value = 64963
$ws.Range("L104").Value = 64963
$ws.Range("N104").Value = -70205
$ws.Range("H122").Value = 6197.5713
$ws.Range("I122").Value = 8992
$ws.Range("K122").Value = 26976
$ws.Range("M122").Value = -24526
$ws.Range("H126").Value = 2999.6667
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -16938.5
$ws.Range("H134").Value = 1573.8718
$ws.Range("I134").Value = 1455.0834
$ws.Range("K134").Value = 4365.2502
$ws.Range("M134").Value = -1830.2502
$ws.Range("H136").Value = 4271.5884
$ws.Range("I136").Value = 3006.6
$ws.Range("J136").Value = 4798.6665
$ws.Range("K136").Value = 9019.799999999999
$ws.Range("L136").Value = 14395.9995
$ws.Range("M136").Value = -6469.799999999999
$ws.Range("N136").Value = -19495.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 501.5
$ws.Range("I69").Value = 637.3333
$ws.Range("J69").Value = 94
$ws.Range("K69").Value = 1911.9999
$ws.Range("L69").Value = 282
$ws.Range("M69").Value = -1100.9999
$ws.Range("N69").Value = -1904
$ws.Range("H72").Value = 501.5
$ws.Range("I72").Value = 637.3333
$ws.Range("J72").Value = 94
$ws.Range("K72").Value = 5735.9997
$ws.Range("L72").Value = 846
$ws.Range("M72").Value = -1679.9997
$ws.Range("N72").Value = -8958
$ws.Range("H109").Value = 4000
$ws.Range("I109").Value = 4000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 12000
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -10960
$ws.Range("H136").Value = 1497.5
$ws.Range("J136").Value = 1495
$ws.Range("L136").Value = 4485
$ws.Range("N136").Value = -14685

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1195.4
$ws.Range("I102").Value = 1137.069
$ws.Range("K102").Value = 1137.069
$ws.Range("M102").Value = 484.931
$ws.Range("H112").Value = 62236.5
$ws.Range("J112").Value = 62236.5
$ws.Range("L112").Value = 62236.5
$ws.Range("N112").Value = -64452.5
$ws.Range("H124").Value = 100122.375
$ws.Range("J124").Value = 111955.8
$ws.Range("L124").Value = 111955.8
$ws.Range("N124").Value = -121775.8
$ws.Range("H126").Value = 2954.8948
$ws.Range("I126").Value = 2078.75
$ws.Range("K126").Value = 6236.25
$ws.Range("M126").Value = -3766.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2975.9644
$ws.Range("I46").Value = 1112
$ws.Range("J46").Value = 3199.64
$ws.Range("K46").Value = 1112
$ws.Range("L46").Value = 3199.64
$ws.Range("M46").Value = -924
$ws.Range("N46").Value = -3575.64
$ws.Range("H55").Value = 836.25
$ws.Range("I55").Value = 963.55554
$ws.Range("J55").Value = 607.1
$ws.Range("K55").Value = 963.55554
$ws.Range("L55").Value = 607.1
$ws.Range("M55").Value = -790.55554
$ws.Range("N55").Value = -953.1
$ws.Range("H59").Value = 46179.832
$ws.Range("J59").Value = 46179.832
$ws.Range("L59").Value = 46179.832
$ws.Range("N59").Value = -47487.832
$ws.Range("H136").Value = 3689.1
$ws.Range("I136").Value = 3693.4119
$ws.Range("K136").Value = 11080.2357
$ws.Range("M136").Value = -8530.235700000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 22790
$ws.Range("I58").Value = 22790
$ws.Range("K58").Value = 22790
$ws.Range("M58").Value = -22482
$ws.Range("H81").Value = 9971.286
$ws.Range("I81").Value = 10099.777
$ws.Range("J81").Value = 9874.916999999999
$ws.Range("K81").Value = 20199.554
$ws.Range("L81").Value = 19749.834
$ws.Range("M81").Value = -19138.554
$ws.Range("N81").Value = -21871.834
$ws.Range("H84").Value = 9971.286
$ws.Range("I84").Value = 10099.777
$ws.Range("J84").Value = 9874.916999999999
$ws.Range("K84").Value = 100997.77
$ws.Range("L84").Value = 98749.17
$ws.Range("M84").Value = -95693.77
$ws.Range("N84").Value = -109357.17
$ws.Range("H114").Value = 44990
$ws.Range("J114").Value = 44990
$ws.Range("L114").Value = 44990
$ws.Range("N114").Value = -53668
$ws.Range("H132").Value = 1704.566
$ws.Range("I132").Value = 1556.8372
$ws.Range("K132").Value = 4670.5116
$ws.Range("M132").Value = -2140.5116
